# "Sweatshirts and fleece" back to "sweatshirts"
#
# Row 6 of the "apparel" sheet currently reads:
#   A6 = "Sweatshirts and Fleece"
#   E6 = "sweatshirts & fleece"
# It should be simplified back to:
#   A6 = "Sweatshirts"
#   E6 = "sweatshirts"
#
# All other data in the sheet is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apparel")

$ws.Range("A6").Value = "Sweatshirts"
$ws.Range("E6").Value = "sweatshirts"

# Leftover cursor/selection position from the editing session.
$ws.Range("E22").Select()
